# Update Selector - (Sharepoint Button file name)
#
# The "Constants" sheet contains a duplicate configuration row: row 28
# (Name = "InputFileName", Value = "Co Code House Bank Account ID") is an
# exact duplicate of row 21. Remove the stray duplicate row so the sheet's
# selector/value list lines up correctly again; everything below shifts
# up by one row as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

$ws.Rows.Item(28).Delete()
